$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teams")

# Insert two new blank columns: C (TeamType) and, after the existing
# Channel1Name shifts to D, a second new column at E (Channel1Type).
# This leaves Channel2Name at F; Channel2Type (G) will just be written
# directly since it's beyond the current used range.
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(5).Insert()

# Fill column-by-column (TeamType, then Channel1Type, then Channel2Type)
# so new shared-string entries are interned in the same order the
# original authoring tool produced them.

# TeamType column
$ws.Range("C1").Value2 = "TeamType"
$ws.Range("C2").Value2 = "Private"
$ws.Range("C3").Value2 = "Private"
$ws.Range("C4").Value2 = "Private"
$ws.Range("C5").Value2 = "Private"

# Channel1Type column
$ws.Range("E1").Value2 = "Channel1Type"
$ws.Range("E2").Value2 = "standard"
$ws.Range("E3").Value2 = "private"
$ws.Range("E4").Value2 = "private"
$ws.Range("E5").Value2 = "private"

# Channel2Type column
$ws.Range("G1").Value2 = "Channel2Type"
$ws.Range("G2").Value2 = "private"
$ws.Range("G3").Value2 = "private"
$ws.Range("G4").Value2 = "standard"
$ws.Range("G5").Value2 = "private"

# Column widths: TeamType narrower, Channel1Type explicit (closest
# achievable character widths to the target stored widths).
$ws.Columns.Item(3).ColumnWidth = 9.3
$ws.Columns.Item(5).ColumnWidth = 13.75

# New selection / active cell on the Teams sheet, then make it the
# active/selected tab (matches activeTab moving from Groups to Teams).
$ws.Range("G2").Select()
$ws.Activate()

$wb.Save()
